$wb = $excel.ActiveWorkbook

# --- Rename sheets (new timestamped names) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912049223886"
$wb.Worksheets.Item(2).Name = "NB_TO-1650291207903396"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912079053903"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912079843876"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912080603914"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912048853917.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912049054203.csv"
$ws1.Range("B4").Value = "go_stims-16502912049063852.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291204921393.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_9-16502912055913858.csv"
$ws2.Range("B3").Value = "OB-16502912065533981.csv"
$ws2.Range("B4").Value = "TB-16502912069513898.csv"
$ws2.Range("B5").Value = "ZB-match_8-16502912054183848.csv"
$ws2.Range("B6").Value = "OB-16502912056764178.csv"
$ws2.Range("B7").Value = "TB-16502912072043881.csv"
$ws2.Range("B8").Value = "ZB-match_0-1650291205623389.csv"
$ws2.Range("B9").Value = "OB-16502912058743901.csv"
$ws2.Range("B10").Value = "TB-16502912078783922.csv"

# --- Sheet 3: RS_TO (no cell value changes, only name) ---

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912079353907.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912079093893.csv"
$ws4.Range("B4").Value = "MM_stims-16502912079663868.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291207936394.csv"
$ws4.Range("B6").Value = "MM_stims-16502912079823866.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912079673882.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502912080463915.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912079883885.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650291208029388.csv"
$ws5.Range("B5").Value = "SAT_stims-16502912080153868.csv"
